$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the duplicate rows: rows 13, 14 (duplicate "John Portman"), and rows 17, 18 (duplicate "Hana Abbass")
# Delete from bottom to top to keep row indices valid
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(17).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(13).Delete()

# Update the selection to match the target state
$ws.Range("A30").Select()
